# Anapa, GL, Gelen. 2019 (added)
#
# This script:
#  1. Inserts a new blank row at position 16, which shifts the existing
#     "Anapa" (2020-2023) and "Gelendzhik" (2020-2023) data blocks down
#     by one row each. Because rows 16/17 (below the Gorjachij Kljuch
#     block) and rows 22/23 (below the now-shifted Anapa block) were
#     already unused/blank, this single insert is enough to open up
#     exactly two new slots per city (2018 placeholder + 2019 data row)
#     without disturbing any existing data further down the sheet.
#  2. Adds the newly available 2018/2019 values for Anapa and
#     Gelendzhik, plus a few previously-missing 2018/2019 figures for
#     Gorjachij Kljuch.
#  3. Restores cell formatting (center alignment == style index "2")
#     on every newly written cell, matching the rest of the table.
#  4. Updates the active-cell selection to match the saved workbook
#     state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

function Set-Cell {
    param(
        [string]$addr,
        $value
    )
    $ws.Range($addr).Value = $value
    $ws.Range($addr).HorizontalAlignment = $xlCenter
}

# ---------------------------------------------------------------
# 1. Make room: insert one blank row before row 16 (before the
#    Anapa block). Everything from old row 16 downward shifts to
#    row+1; this also carries the Gelendzhik block down by one row.
# ---------------------------------------------------------------
$ws.Rows.Item(16).Insert()

# Excel's row insert copies the formatting (incl. style) of the row
# above into every column of the freshly inserted row. Only columns
# A and B of this new row actually hold data (2018 placeholder), so
# drop the inherited formatting/cells from the rest of the columns.
$ws.Range("C16:U16").Clear()

# ---------------------------------------------------------------
# 2. Gorjachij Kljuch (Горячий ключ) - fill newly reported figures
# ---------------------------------------------------------------

# Row 10 - 2018: add saldo (U)
Set-Cell "U10" 1872

# Row 11 - 2019: add avgemployers(D), unemployed(E), companies(N),
# conscap(P), consnewareas(Q), consnewapt(R), saldo(U)
Set-Cell "D11" 8.4
Set-Cell "E11" 267
Set-Cell "N11" 942
Set-Cell "P11" 94.8
Set-Cell "Q11" 73.932
Set-Cell "R11" 1157
Set-Cell "U11" 1717

# ---------------------------------------------------------------
# 3. Anapa (Анапа) - new rows 16 (2018) and 17 (2019)
# ---------------------------------------------------------------

# Row 16 - 2018 (only name/year known so far)
Set-Cell "A16" "Анапа"
Set-Cell "B16" 2018

# Row 17 - 2019 (full data)
Set-Cell "A17" "Анапа"
Set-Cell "B17" 2019
Set-Cell "D17" 26
Set-Cell "E17" 743
Set-Cell "F17" 33185
Set-Cell "N17" 3822
Set-Cell "O17" 1871.8
Set-Cell "P17" 3351.1
Set-Cell "Q17" 484.169
Set-Cell "R17" 9054
Set-Cell "S17" 23432.8
Set-Cell "T17" 749.8
Set-Cell "U17" 10533

# ---------------------------------------------------------------
# 4. Gelendzhik (Геленджик) - new rows 22 (2018) and 23 (2019)
# ---------------------------------------------------------------

# Row 22 - 2018 (only name/year known so far)
Set-Cell "A22" "Геленджик"
Set-Cell "B22" 2018

# Row 23 - 2019 (full data)
Set-Cell "A23" "Геленджик"
Set-Cell "B23" 2019
Set-Cell "D23" 17.2
Set-Cell "E23" 266
Set-Cell "F23" 36573
Set-Cell "N23" 2717
Set-Cell "O23" 2218.7
Set-Cell "P23" 200.3
Set-Cell "Q23" 85.383
Set-Cell "R23" 915
Set-Cell "S23" 14798.4
Set-Cell "T23" 878.2
Set-Cell "U23" -553

# ---------------------------------------------------------------
# 5. Final selection, to match saved workbook view state
# ---------------------------------------------------------------
$ws.Range("P28").Select()
